$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column R: "Update Only" header with "No" for every data row
# (sharedStrings gains "Update Only" / "No"; dimension grows to A1:R9).
$ws.Range("R1").Value = "Update Only"

for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 18).Value = "No"
}

# Match the author's final view state: window scrolled so column D is
# leftmost, and R3:R9 selected (was L1:L9 before the edit).
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$ws.Range("R3:R9").Select()
